$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("G2").Value = "2016-08-18 11:07:58"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("H2").Value = "2016-08-18 11:07:52"
$ws2.Range("K2").Value = "2016-08-18 11:08:16"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("K2").Value = "2016-08-18 11:08:24"
